$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in existing date string (10/16 -> 10/17)
$ws.Range("B22").Value = "10/17 /1:00"

# Row 23: new meeting date/time entry and attendance marks
$ws.Range("B23").Value = "10/20 /4:15"
$ws.Range("D23").Value = "A"
$ws.Range("E23").Value = "A"
$ws.Range("F23").Value = "A"
$ws.Range("G23").Value = "A"
$ws.Range("H23").Value = "A"
$ws.Range("I23").Value = "A"

# Update view state (scroll position + active selection)
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("I24").Select()
